$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Clear the old C1:D1 leftovers from the prior 1x4 table layout; the new
# table is a 2-column (A:B) lookup of local certificate image paths to the
# remote certificate URLs used to build it.
$ws2.Range("C1:D1").ClearContents()

# Seed the new local-image-path column first (so the new shared strings are
# interned in this order), then fill in the header row, then the URL column
# that reuses the pre-existing shared strings.
$ws2.Range("A2").Value = "../static/images/Certifications/Certificate-Python-Bootcamp.png"
$ws2.Range("A3").Value = "../static/images/Certifications/Certificate-Python-Django-Bootcamp.png"
$ws2.Range("A5").Value = "../static/images/Certifications/PR Certification.png"
$ws2.Range("A4").Value = "../static/images/Certifications/SEO Certification.png"

$ws2.Range("A1").Value = "img"
$ws2.Range("B1").Value = "cert"

$ws2.Range("B2").Value = "https://udemy-certificate.s3.amazonaws.com/image/UC-b527ed3c-4fd0-4701-bc19-09d06acf6bfe.jpg"
$ws2.Range("B3").Value = "https://udemy-certificate.s3.amazonaws.com/image/UC-c2d0c32d-01d5-4c24-9c77-ff7e6ace6c2f.jpg"
$ws2.Range("B4").Value = "https://udemy-certificate.s3.amazonaws.com/image/UC-17644286-b7e4-4c68-a168-b5aaabfa8596.jpg"
$ws2.Range("B5").Value = "https://udemy-certificate.s3.amazonaws.com/image/UC-5e882f16-f8d8-4885-a604-b41d60407052.jpg"

# Widen column A so the long relative image paths are readable.
$ws2.Columns.Item(1).ColumnWidth = 12.6

# New table needs a print layout just like Sheet1 already has.
$ws2.PageSetup.Orientation = 1

# Bring Sheet2 to the front (moves tabSelected/activeTab off of Sheet1) and
# leave the cursor parked where the author left off, a few rows under the
# new table.
$ws2.Activate()
$ws2.Range("E8").Select()
